$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("全件")

$ws.Range("C4").Value = "上記期間は無休"

$ws.Range("C19").Value = "8月11日(月)"

$ws.Range("C21").Value = "8月17日(日)"
$ws.Range("E21").Value = 18

$ws.Range("C48").Value = "8月15日(金), 8月16日(土), 8月17日(日)"
$ws.Range("E48").Value = "19に配達"

$ws.Range("C49").Value = "8月15日(金), 8月16日(土), 8月17日(日)"

$ws.Range("C53").Value = "8月10日(日), 8月11日(月), 8月14日(木), 8月15日(金), 8月16日(土), 8月17日(日)"

$ws.Range("C56").Value = "8月15日(金)"

$ws.Range("C76").Value = "上記期間は全て休業"

$ws.Range("C82").Value = "8月10日(日), 8月11日(月), 8月12日(火), 8月13日(水), 8月14日(木), 8月15日(金), 8月16日(土), 8月17日(日)"

$ws.Range("C95").Value = "8月11日(月)"
$ws.Range("E95").Value = "15　未定"

$ws.Range("C98").Value = "8月10日(日), 8月17日(日)"
$ws.Range("E98").Value = "定休日、日曜日"

$ws.Range("C106").Value = "8月13日(水), 8月14日(木), 8月15日(金)"

$ws.Range("C111").Value = "8月10日(日), 8月11日(月), 8月17日(日)"

$ws.Range("C114").Value = "8月12日(火)"

$ws.Range("C115").Value = "8月11日(月)"

$ws.Range("C116").Value = "8月11日(月)"

$ws.Range("C122").Value = "8月10日(日), 8月15日(金), 8月17日(日)"

$ws.Range("C127").Value = "8月10日(日), 8月11日(月), 8月12日(火), 8月13日(水), 8月14日(木)"

$ws.Range("C137").Value = "8月15日(金)"

$ws.Range("C147").Value = "8月10日(日), 8月17日(日)"

$ws.Range("C152").Value = "8月10日(日), 8月17日(日)"

$ws.Range("C160").Value = "上記期間は無休"

$ws.Range("C167").Value = "8月10日(日), 8月17日(日)"

$ws.Range("C168").Value = "8月15日(金)"

$ws.Range("C180").Value = "8月15日(金), 8月16日(土), 8月17日(日)"

$ws.Range("C204").Value = "8月15日(金), 8月17日(日)"

$ws.Range("C207").Value = "8月10日(日), 8月17日(日)"

$ws.Range("C209").Value = "8月15日(金), 8月17日(日)"

$ws.Range("C220").Value = "上記期間は無休"
$ws.Range("D220").Value = "120人"
